$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the cell contents: A1 becomes numeric 0, C1 becomes the text "s"
$ws.Range("A1").Value = 0
$ws.Range("C1").Value = "s"

# Move the active selection to R15
$ws.Range("R15").Select()
